# fixed ubuntu ami-ids for aws
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 holds the "Ubuntu 14" AMI ids per region. Refresh them, and add
# the previously-missing EU (London) AMI in column J.
$ws.Range("D4").Value = "ami-9dde7f8b"
$ws.Range("F4").Value = "ami-9d772efd"
$ws.Range("G4").Value = "ami-0e2aa66e"
$ws.Range("H4").Value = "ami-115d7777"
$ws.Range("I4").Value = "ami-6039ed0f"
$ws.Range("J4").Value = "ami-c29184a6"

# Selection moved to J7 as part of the edit.
$ws.Range("J7").Select()
